$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) is treated as text so values such as
# "604.23" or "70.644.43" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "70.644.43"
$ws.Range("E2").Value = "  +0.78%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.622.83"
$ws.Range("E3").Value = "  +2.31%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "604.23"
$ws.Range("E5").Value = "  +0.02%  "

# Row 6 - Solana
$ws.Range("D6").Value = "196.62"
$ws.Range("E6").Value = "  -0.14%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  +0.22%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.04%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.94%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  -1.11%  "

# Row 11 - Avalanche
$ws.Range("D11").Value = "53.71"
$ws.Range("E11").Value = "  -0.65%  "

# Row 12 - ShibaInu
$ws.Range("E12").Value = "  +0.43%  "

# Row 13 - Polkadot
$ws.Range("E13").Value = "  +0.27%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.195.55"
$ws.Range("E14").Value = "  +2.13%  "

# Row 15 - was Uniswap, now BitcoinCash
$ws.Range("B15").Value = "BitcoinCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D15").Value = "601.49"
$ws.Range("E15").Value = "  -0.38%  "

# Row 16 - was BitcoinCash, now Uniswap
$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").Value = "13.03"
$ws.Range("E16").Value = "  +2.42%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "70.703.99"
$ws.Range("E17").Value = "  +0.71%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.640.56"
$ws.Range("E18").Value = "  +2.93%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "19.11"
$ws.Range("E19").Value = "  -0.67%  "

# Row 21 - Polygon
$ws.Range("D21").Value = "0.998"
$ws.Range("E21").Value = "  +0.30%  "

# Row 22 - InternetComputer(DFINITY)
$ws.Range("D22").Value = "17.84"
$ws.Range("E22").Value = "  -0.85%  "

# Row 23 - Toncoin
$ws.Range("D23").Value = "5.19"
$ws.Range("E23").Value = "  -1.52%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "101.76"
$ws.Range("E24").Value = "  -1.02%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +0.43%  "

# Row 26 - ImmutableX
$ws.Range("E26").Value = "  -3.57%  "

# Row 27 - RenderToken
$ws.Range("D27").Value = "10.78"
$ws.Range("E27").Value = "  -1.84%  "

# Row 28 - Filecoin
$ws.Range("E28").Value = "  +0.02%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "33.92"
$ws.Range("E29").Value = "  +0.38%  "

# Row 30 - dogwifhat
$ws.Range("D30").Value = "4.66"
$ws.Range("E30").Value = "  +6.36%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "7.24"
$ws.Range("E31").Value = "  +1.41%  "

# Row 32 - Cosmos
$ws.Range("E32").Value = "  -2.64%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  +1.72%  "

# Row 34 - OKB
$ws.Range("D34").Value = "63.60"
$ws.Range("E34").Value = "  +0.31%  "

# Row 35 - PEPE
$ws.Range("D35").Value = "0.0₃0898"
$ws.Range("E35").Value = "  +6.95%  "

# Row 36 - Maker
$ws.Range("D36").Value = "3.918.35"
$ws.Range("E36").Value = "  +3.72%  "

# Row 37 - Bittensor
$ws.Range("D37").Value = "544.77"
$ws.Range("E37").Value = "  +10.98%  "

# Row 38 - Fetch.AI
$ws.Range("D38").Value = "3.12"
$ws.Range("E38").Value = "  +1.54%  "

# Row 39 - Dai
$ws.Range("E39").Value = "  +0.11%  "

# Row 40 - InjectiveProtocol
$ws.Range("D40").Value = "37.07"
$ws.Range("E40").Value = "  +0.86%  "

# Row 41 - TheGraph
$ws.Range("E41").Value = "  -1.21%  "

# Row 42 - Stacks
$ws.Range("D42").Value = "3.54"
$ws.Range("E42").Value = "  -3.29%  "

# Row 43 - Kaspa
$ws.Range("E43").Value = "  -0.33%  "

# Row 45 - ApeXProtocol
$ws.Range("D45").Value = "3.40"
$ws.Range("E45").Value = "  +3.10%  "

# Row 46 - ThetaToken
$ws.Range("D46").Value = "2.87"
$ws.Range("E46").Value = "  +0.71%  "

# Row 47 - Stellar
$ws.Range("E47").Value = "  +0.27%  "

# Row 48 - THORChain
$ws.Range("D48").Value = "8.63"
$ws.Range("E48").Value = "  -0.62%  "

# Row 49 - FirstDigitalUSD
$ws.Range("E49").Value = "  -0.26%  "

# Row 50 - FLOKI
$ws.Range("E50").Value = "  +0.93%  "

# Row 51 - Mantle
$ws.Range("E51").Value = "  +2.69%  "

# Restore the default (Normal) style on column D so no stray cell-level
# style attributes are introduced by the temporary text formatting above.
$ws.Range("D2:D51").Style = "Normal"
